# Auto-generated Excel COM-interop script
# Applies the "Updated cryptos list" data refresh:
#  - Updates Price (column D) and Volume(1h) (column E) values for many rows
#  - Swaps row order data for Litecoin/Fetch.AI (rows 25-26)
#  - Swaps row order data for PancakeSwap/ImmutableX (rows 30-31)
#
# Cells are written as text (NumberFormat "@") and then reset to the
# "Normal" style so that numeric-looking strings such as "592.40" or
# "0.530" are preserved verbatim instead of being parsed into floating
# point numbers (which would silently drop significant trailing zeros
# or punctuation used as thousands separators).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}


Set-TextValue "D2" '63.826.69'
Set-TextValue "E2" '  -0.23%  '
Set-TextValue "D3" '3.144.57'
Set-TextValue "E3" '  +0.34%  '
Set-TextValue "E4" '  +0.09%  '
Set-TextValue "D5" '592.40'
Set-TextValue "E5" '  +0.27%  '
Set-TextValue "D6" '145.37'
Set-TextValue "E6" '  -1.30%  '
Set-TextValue "E7" '  +0.07%  '
Set-TextValue "D8" '3.135.61'
Set-TextValue "E8" '  +0.19%  '
Set-TextValue "D9" '0.530'
Set-TextValue "E9" '  -0.77%  '
Set-TextValue "E10" '  -0.22%  '
Set-TextValue "D11" '5.87'
Set-TextValue "E11" '  +2.08%  '
Set-TextValue "E12" '  -2.10%  '
Set-TextValue "E13" '  -2.76%  '
Set-TextValue "D14" '37.23'
Set-TextValue "E14" '  -0.70%  '
Set-TextValue "D15" '3.665.36'
Set-TextValue "E15" '  +0.32%  '
Set-TextValue "E16" '  -1.41%  '
Set-TextValue "E17" '  +1.80%  '
Set-TextValue "D18" '3.141.37'
Set-TextValue "E18" '  +0.07%  '
Set-TextValue "D19" '63.719.69'
Set-TextValue "E19" '  -0.19%  '
Set-TextValue "D20" '468.63'
Set-TextValue "E20" '  +0.00%  '
Set-TextValue "E21" '  -0.15%  '
Set-TextValue "E22" '  -0.49%  '
Set-TextValue "E23" '  -0.64%  '
Set-TextValue "D24" '13.00'
Set-TextValue "E24" '  -2.34%  '
Set-TextValue "B25" 'Fetch.AI'
Set-TextValue "C25" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D25" '2.33'
Set-TextValue "E25" '  +6.31%  '
Set-TextValue "B26" 'Litecoin'
Set-TextValue "C26" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D26" '81.36'
Set-TextValue "E26" '  -1.43%  '
Set-TextValue "E27" '  +0.04%  '
Set-TextValue "D28" '9.80'
Set-TextValue "E28" '  +8.89%  '
Set-TextValue "D29" '7.40'
Set-TextValue "E29" '  +7.83%  '
Set-TextValue "B30" 'ImmutableX'
Set-TextValue "C30" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D30" '2.24'
Set-TextValue "E30" '  +0.09%  '
Set-TextValue "B31" 'PancakeSwap'
Set-TextValue "C31" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D31" '2.71'
Set-TextValue "E31" '  -0.44%  '
Set-TextValue "E32" '  +0.11%  '
Set-TextValue "D33" '27.79'
Set-TextValue "E33" '  +2.15%  '
Set-TextValue "E34" '  +1.94%  '
Set-TextValue "D35" '0.0₃0844'
Set-TextValue "E35" '  -4.17%  '
Set-TextValue "E36" '  +1.02%  '
Set-TextValue "E37" '  -3.33%  '
Set-TextValue "D38" '6.16'
Set-TextValue "E38" '  +0.32%  '
Set-TextValue "D39" '3.23'
Set-TextValue "E39" '  -5.31%  '
Set-TextValue "D40" '51.47'
Set-TextValue "E40" '  +0.87%  '
Set-TextValue "D41" '9.29'
Set-TextValue "E41" '  +6.34%  '
Set-TextValue "D42" '455.29'
Set-TextValue "E42" '  -0.16%  '
Set-TextValue "E43" '  +5.32%  '
Set-TextValue "D44" '0.0372'
Set-TextValue "E44" '  -0.37%  '
Set-TextValue "D45" '2.911.77'
Set-TextValue "E45" '  +0.28%  '
Set-TextValue "D46" '39.73'
Set-TextValue "E46" '  +11.48%  '
Set-TextValue "E47" '  -3.26%  '
Set-TextValue "D48" '130.24'
Set-TextValue "E48" '  +2.84%  '
Set-TextValue "D50" '2.25'
Set-TextValue "E50" '  +2.59%  '
Set-TextValue "D51" '0.111'
Set-TextValue "E51" '  -1.07%  '
